$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear existing data + hyperlinks so we can rewrite the full table cleanly
$ws.Cells.Clear()
$ws.Hyperlinks.Delete()

# Header row
$ws.Range("A1").Value = '取得日時'
$ws.Range("B1").Value = 'タイトル'
$ws.Range("C1").Value = 'カテゴリ'
$ws.Range("D1").Value = '価格'
$ws.Range("E1").Value = '締切'
$ws.Range("F1").Value = 'URL'
$ws.Range("G1").Value = '優先度スコア'
$ws.Range("H1").Value = 'スキル概要'

# Data rows
# Row 2
$ws.Range("A2").Value = '2025-11-13 12:39:17'
$ws.Range("B2").Value = 'AIエンジニア募集|LLM・LangChain・RAG・Python経験者歓迎'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5433318')
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("G2").Value = 505
$ws.Range("H2").Value = '🔥AI,Python'

# Row 3
$ws.Range("A3").Value = '2025-11-13 12:39:17'
$ws.Range("B3").Value = '大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5427956')
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("G3").Value = 310
$ws.Range("H3").Value = '🔥AI,Ai'

# Row 4
$ws.Range("A4").Value = '2025-11-13 12:39:17'
$ws.Range("B4").Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5217096')
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("G4").Value = 243
$ws.Range("H4").Value = '🔥API ◆ツール'

# Row 5
$ws.Range("A5").Value = '2025-11-13 12:39:17'
$ws.Range("B5").Value = '【スポット】画像&動画 編集・公開サービスサイトのTypeScript,Node.JSでの更新作業'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5433199')
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("G5").Value = 190
$ws.Range("H5").Value = '🔥TypeScript ◆Node.js ◇サイト'

# Row 6
$ws.Range("A6").Value = '2025-11-13 12:39:17'
$ws.Range("B6").Value = '【急募】九九アプリの開発依頼'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5433544')
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("G6").Value = 85
$ws.Range("H6").Value = '◆開発 ◇アプリ'

# Row 7
$ws.Range("A7").Value = '2025-11-13 12:39:17'
$ws.Range("B7").Value = '【日本人限定・長期募集】SNS運用担当募集|Web開発会社 JapanDream'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5432819')
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("G7").Value = 68
$ws.Range("H7").Value = '◆開発'

# Row 8
$ws.Range("A8").Value = '2025-11-13 12:39:17'
$ws.Range("B8").Value = '【急募】オンラインガチャ制作のフリーランスを探しています!'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5433143')
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("G8").Value = 25

# Row 9
$ws.Range("A9").Value = '2025-11-13 12:39:17'
$ws.Range("B9").Value = 'Networkエンジニア'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5432661')
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("G9").Value = 25

# Row 10
$ws.Range("A10").Value = '2025-11-13 12:39:17'
$ws.Range("B10").Value = 'WooCommerce消費税設定のカスタマイズ依頼'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5432929')
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("G10").Value = 18

# Row 11
$ws.Range("A11").Value = '2025-11-13 12:39:17'
$ws.Range("B11").Value = '適合商品検索ページ作成'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5432621')
$ws.Range("F11").Style = "Hyperlink"
$ws.Range("G11").Value = 18

# Row 12
$ws.Range("A12").Value = '2025-11-13 12:39:17'
$ws.Range("B12").Value = 'アカウント運用代行募集 経験者のみ 無在庫販売 EC販売をご経験のある方のみで長期募集'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5433373')
$ws.Range("F12").Style = "Hyperlink"
$ws.Range("G12").Value = 10

# Column widths
$ws.Columns.Item(2).ColumnWidth = 50.166666666666664
$ws.Columns.Item(8).ColumnWidth = 26.166666666666668
